$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data re-pull for the 2023-01-09 build replaced a handful of placeholder
# values in the syllabus sheet with the actual professor / date strings.

$carlos = "519033 - Carlos Yujiro Shigue"
$emerson = "7290967 - Emerson Gonçalves de Melo"

# Row 10 - "Objetivos:" value cells
$ws.Range("B10").Value = $carlos
$ws.Range("C10").Value = $carlos

# Row 13 - "Programa resumido:" value cells (date-looking text; write it via
# a quoted formula then paste-special as values so Excel keeps it as a
# plain text string instead of auto-converting it to a date serial, while
# preserving the existing cell styles).
$ws.Range("B13").Formula = "=""01/01/2023"""
$ws.Range("B13").Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("C13").Formula = "=""01/01/2023"""
$ws.Range("C13").Copy()
$ws.Range("C13").PasteSpecial(-4163)

# Row 15 - "Programa:" value cells
$ws.Range("B15").Value = $carlos
$ws.Range("C15").Value = $carlos

# Row 18 - "Método:" value cells
$ws.Range("B18").Value = $emerson
$ws.Range("C18").Value = $emerson
